$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the SHHConfig row (row 7) with the new configuration file path,
# replacing the old "ssh" placeholder value across all experiment columns (B:AW).
$ws.Range("B7:AW7").Value = "Configurations/SSHConfig.csv"

# Reflect the post-edit selection/view state: the user had just filled the row
# and ended up with B7:AW7 selected, scrolled so column AW is the leftmost visible.
$ws.Range("B7:AW7").Select()
$excel.ActiveWindow.ScrollColumn = $ws.Range("AW1").Column

$wb.Save()
